$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the existing value at F215 (recomputed FX_EURUSD figure)
$ws.Range("F215").Value = 1.19657301902771

# New rows of portfolio data (216-220)
$data = @(
    @(46055, 2110.53, 114482.6636669922, 114482.6636669922, 6976.43994140625, 1.184918284416199),
    @(46056, 2110.53, 114437.6,          114437.6,          6917.81005859375, 1.179871439933777),
    @(46057, 2110.53, 114107.16,         114107.16,         6882.72021484375, 1.18161404132843),
    @(46058, 2110.53, 113371.18,         113371.18,         6798.39990234375, 1.180163860321045),
    @(46059, 2110.53, 116172.4136669922, 116172.4136669922, 6932.2998046875,  1.177786946296692)
)

# Copy the date cell formatting (border/font/alignment/number format) from
# the last existing data row (A215) down onto the new date cells.
$ws.Range("A215").Copy()
$ws.Range("A216:A220").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$row = 216
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}
